$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.0292345
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.4428765120700495
$ws.Range("J2").Value = 0.346386487911515
$ws.Range("M2").Value = 0.071358
$ws.Range("N2").Value = 0.142716
$ws.Range("O2").Value = 0.01919591193090569
$ws.Range("P2").Value = 0.01411929935366186
$ws.Range("Q2").Value = 0.002086115451
$ws.Range("R2").Value = 0.008344461804000001
$ws.Range("S2").Value = 0.008501418521963362
$ws.Range("T2").Value = 0.004890734514886253
# Row 3
$ws.Range("G3").Value = 0.0292345
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.4428765120700495
$ws.Range("J3").Value = 0.346386487911515
$ws.Range("O3").Value = 0.7144019644080171
$ws.Range("P3").Value = 0.7882033865305114
$ws.Range("Q3").Value = 0.0776376231325
$ws.Range("R3").Value = 0.465825738795
$ws.Range("S3").Value = 0.3163918502130142
$ws.Range("T3").Value = 0.2730230028202661
# Row 4
$ws.Range("G4").Value = 0.0292345
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.4428765120700495
$ws.Range("J4").Value = 0.346386487911515
$ws.Range("M4").Value = 0.9728370000000001
$ws.Range("N4").Value = 1.945674
$ws.Range("O4").Value = 0.2617014683024538
$ws.Range("P4").Value = 0.1924910567184946
$ws.Range("Q4").Value = 0.0284404032765
$ws.Range("R4").Value = 0.113761613106
$ws.Range("S4").Value = 0.1159014334854014
$ws.Range("T4").Value = 0.06667630109109558
# Row 5
$ws.Range("G5").Value = 0.0292345
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.4428765120700495
$ws.Range("J5").Value = 0.346386487911515
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.017474
$ws.Range("N5").Value = 0.052422
$ws.Range("O5").Value = 0.004700655358623364
$ws.Range("P5").Value = 0.005186257397332197
$ws.Range("Q5").Value = 0.000510843653
$ws.Range("R5").Value = 0.003065061918
$ws.Range("S5").Value = 0.002081809849670503
$ws.Range("T5").Value = 0.001796449485267014
# Row 6
$ws.Range("I6").Value = 0.5571234879299505
$ws.Range("J6").Value = 0.6536135120884849
$ws.Range("M6").Value = 0.071358
$ws.Range("N6").Value = 0.142716
$ws.Range("O6").Value = 0.01919591193090569
$ws.Range("P6").Value = 0.01411929935366186
$ws.Range("Q6").Value = 0.002624261808
$ws.Range("R6").Value = 0.015745570848
$ws.Range("S6").Value = 0.01069449340894233
$ws.Range("T6").Value = 0.009228564838775601
# Row 7
$ws.Range("I7").Value = 0.5571234879299505
$ws.Range("J7").Value = 0.6536135120884849
$ws.Range("O7").Value = 0.7144019644080171
$ws.Range("P7").Value = 0.7882033865305114
$ws.Range("S7").Value = 0.3980101141950028
$ws.Range("T7").Value = 0.5151803837102452
# Row 8
$ws.Range("I8").Value = 0.5571234879299505
$ws.Range("J8").Value = 0.6536135120884849
$ws.Range("M8").Value = 0.9728370000000001
$ws.Range("N8").Value = 1.945674
$ws.Range("O8").Value = 0.2617014683024538
$ws.Range("P8").Value = 0.1924910567184946
$ws.Range("Q8").Value = 0.035777053512
$ws.Range("R8").Value = 0.214662321072
$ws.Range("S8").Value = 0.1458000348170524
$ws.Range("T8").Value = 0.125814755627399
# Row 9
$ws.Range("I9").Value = 0.5571234879299505
$ws.Range("J9").Value = 0.6536135120884849
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.017474
$ws.Range("N9").Value = 0.052422
$ws.Range("O9").Value = 0.004700655358623364
$ws.Range("P9").Value = 0.005186257397332197
$ws.Range("Q9").Value = 0.0006426238239999999
$ws.Range("R9").Value = 0.005783614416
$ws.Range("S9").Value = 0.002618845508952861
$ws.Range("T9").Value = 0.003389807912065182
